# Adds "rename" localization context to the node tree:
#  - Two new rows inserted before the old row 60 ("Left"/左):
#      row 60: InternalRenameNode / Rename / 名称の変更
#      row 61: RenameNode         / Rename / 名称の変更
#  - A new trailing row appended for the Cancel button:
#      row 168: Cancel / Cancel / キャンセル

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 60 ("Left"), pushing the
# rest of the table down. The new rows pick up the surrounding row
# format (same as every other data row in the sheet).
$ws.Range("A60:A61").EntireRow.Insert()

$ws.Cells.Item(60, 1).Value = "InternalRenameNode"
$ws.Cells.Item(60, 2).Value = "Rename"
$ws.Cells.Item(60, 3).Value = "名称の変更"

$ws.Cells.Item(61, 1).Value = "RenameNode"
$ws.Cells.Item(61, 2).Value = "Rename"
$ws.Cells.Item(61, 3).Value = "名称の変更"

# Append the new Cancel row at the end of the table (row 168).
$ws.Cells.Item(168, 1).Value = "Cancel"
$ws.Cells.Item(168, 2).Value = "Cancel"
$ws.Cells.Item(168, 3).Value = "キャンセル"

# Match the author's final selection (the freshly added row).
$ws.Range("A168:C168").Select()
